# Update "想去人数" (F column) counters across the four sheets to the
# values captured at the time the gh-pages output was regenerated
# (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2908
$ws.Range("F3").Value = 21356
$ws.Range("F4").Value = 105
$ws.Range("F5").Value = 3218
$ws.Range("F6").Value = 832
$ws.Range("F7").Value = 621
$ws.Range("F8").Value = 538
$ws.Range("F9").Value = 790
$ws.Range("F10").Value = 296
$ws.Range("F14").Value = 547
$ws.Range("F15").Value = 188
$ws.Range("F16").Value = 307
$ws.Range("F17").Value = 31
$ws.Range("F18").Value = 439
$ws.Range("F19").Value = 117
$ws.Range("F21").Value = 29
$ws.Range("F22").Value = 53
$ws.Range("F23").Value = 140

# --- 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 5

# --- 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6155
$ws.Range("F4").Value = 717
$ws.Range("F5").Value = 1691
$ws.Range("F6").Value = 71

# --- 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6155
$ws.Range("F4").Value = 717
$ws.Range("F5").Value = 1691
$ws.Range("F6").Value = 2908
$ws.Range("F7").Value = 21356
$ws.Range("F12").Value = 3218
$ws.Range("F13").Value = 832
$ws.Range("F15").Value = 71
$ws.Range("F16").Value = 621
$ws.Range("F17").Value = 538
$ws.Range("F18").Value = 790
$ws.Range("F19").Value = 296
$ws.Range("F28").Value = 547
$ws.Range("F30").Value = 188
$ws.Range("F32").Value = 307
$ws.Range("F35").Value = 31
$ws.Range("F36").Value = 439
$ws.Range("F38").Value = 117
$ws.Range("F42").Value = 29
$ws.Range("F43").Value = 54
$ws.Range("F44").Value = 5
$ws.Range("F49").Value = 140
